$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.156.57"
$ws.Range("E2").Value = "  +2.90%  "
$ws.Range("D3").Value = "2.062.80"
$ws.Range("E3").Value = "  +2.32%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'230.21"
$ws.Range("E5").Value = "  +1.70%  "
$ws.Range("D6").Value = "'0.616"
$ws.Range("E6").Value = "  +1.48%  "
$ws.Range("D7").Value = "'58.20"
$ws.Range("E7").Value = "  +6.37%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  +2.20%  "
$ws.Range("D10").Value = "'0.0807"
$ws.Range("E10").Value = "  +2.28%  "
$ws.Range("E11").Value = "  -1.19%  "
$ws.Range("D12").Value = "2.366.80"
$ws.Range("E12").Value = "  +2.35%  "
$ws.Range("D13").Value = "'14.64"
$ws.Range("E13").Value = "  +2.80%  "
$ws.Range("E14").Value = "  +2.39%  "
$ws.Range("D15").Value = "'0.754"
$ws.Range("E15").Value = "  +1.74%  "
$ws.Range("E16").Value = "  +3.20%  "
$ws.Range("D17").Value = "2.062.80"
$ws.Range("E17").Value = "  +2.18%  "
$ws.Range("D18").Value = "38.024.14"
$ws.Range("E18").Value = "  +2.92%  "
$ws.Range("E19").Value = "  +2.22%  "
$ws.Range("D20").Value = "'69.86"
$ws.Range("E20").Value = "  +1.54%  "
$ws.Range("E21").Value = "  +1.61%  "
$ws.Range("D22").Value = "'224.79"
$ws.Range("E22").Value = "  +0.47%  "
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("E24").Value = "  +0.76%  "
$ws.Range("E25").Value = "  +3.19%  "
$ws.Range("D26").Value = "'9.30"
$ws.Range("E26").Value = "  +1.18%  "
$ws.Range("D27").Value = "'165.70"
$ws.Range("E27").Value = "  +0.48%  "
$ws.Range("D28").Value = "'0.135"
$ws.Range("E28").Value = "  +7.84%  "
$ws.Range("E29").Value = "  +1.68%  "
$ws.Range("D30").Value = "'19.06"
$ws.Range("E30").Value = "  +1.90%  "
$ws.Range("E31").Value = "  +1.93%  "
$ws.Range("E32").Value = "  +1.16%  "
$ws.Range("E33").Value = "  +4.16%  "
$ws.Range("E34").Value = "  +0.31%  "
$ws.Range("E35").Value = "  +7.16%  "
$ws.Range("E36").Value = "  +1.57%  "
$ws.Range("D37").Value = "'6.05"
$ws.Range("E37").Value = "  +12.48%  "
$ws.Range("D38").Value = "'3.31"
$ws.Range("E38").Value = "  +5.61%  "
$ws.Range("D39").Value = "'1.00"
$ws.Range("E39").Value = "  -0.36%  "
$ws.Range("D40").Value = "'98.46"
$ws.Range("E40").Value = "  +3.57%  "
$ws.Range("E41").Value = "  +0.80%  "
$ws.Range("D42").Value = "1.478.81"
$ws.Range("E42").Value = "  +0.33%  "
$ws.Range("D43").Value = "'0.0946"
$ws.Range("E43").Value = "  +2.74%  "
$ws.Range("E44").Value = "  +4.19%  "
$ws.Range("D45").Value = "'16.80"
$ws.Range("E45").Value = "  +1.80%  "
$ws.Range("E46").Value = "  +0.06%  "
$ws.Range("E47").Value = "  +17.29%  "
$ws.Range("E48").Value = "  +1.24%  "
$ws.Range("E49").Value = "  +2.04%  "
$ws.Range("E50").Value = "  -1.55%  "
$ws.Range("D51").Value = "2.254.53"
$ws.Range("E51").Value = "  +2.27%  "
